$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Hommes, JEL,2019/20"
$ws.Range("B12").Value = "mon pol should be more aggressive on things to add negative feedback, makes the system more stable"

$ws.Range("A13").Value = "Gabaix 2019, bounded rational NK"
$ws.Range("B13").Value = "opt mon pol isn't price level targeting - truly, the more BR firms are, the less forward-looking, and the less commitment the CB can implement."

$ws.Range("B12:B13").WrapText = $true

$ws.Range("B14").Select()
